# "reset position for all scene"
# The Scene sheet lists one scene configuration per row (rows 11-16).
# Column G holds "RelivePos" (the scene's revive/reset position, a
# "x,y,z" string). Every scene's RelivePos is reset to a new coordinate.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G11").Value = "89,104,0"
$ws.Range("G12").Value = "89,102,0"
$ws.Range("G13").Value = "89,104,0"
$ws.Range("G14").Value = "100,115,0"
$ws.Range("G15").Value = "95,102,0"
$ws.Range("G16").Value = "88,105,0"

# Match the author's final selection/view state as closely as this
# runtime allows (active cell moves to G16).
$ws.Range("G16").Select()
